$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 245
$ws1.Range("F9").Value  = 6781
$ws1.Range("F15").Value = 1091
$ws1.Range("F16").Value = 16159
$ws1.Range("F19").Value = 329
$ws1.Range("F20").Value = 179
$ws1.Range("F22").Value = 11333
$ws1.Range("F24").Value = 967
$ws1.Range("F25").Value = 4460
$ws1.Range("F26").Value = 310
$ws1.Range("F29").Value = 662
$ws1.Range("F30").Value = 319
$ws1.Range("F31").Value = 140

# Sheet 4 = "全部类型" (All Types) - mirrors sheet 1's rows
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 245
$ws4.Range("F10").Value = 6781
$ws4.Range("F17").Value = 1091
$ws4.Range("F18").Value = 16159
$ws4.Range("F21").Value = 329
$ws4.Range("F22").Value = 179
$ws4.Range("F26").Value = 11333
$ws4.Range("F28").Value = 967
$ws4.Range("F29").Value = 4460
$ws4.Range("F30").Value = 310
$ws4.Range("F33").Value = 664
$ws4.Range("F34").Value = 319
$ws4.Range("F35").Value = 140
